$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated quantities (F), computed values (G), and recalculated
# sub-totals / grand totals (B) as captured in the source diff.
$updates = @{
    "F11" = 19
    "G11" = 7375.23
    "B18" = 45217.55
    "F25" = 45
    "G25" = 1268.1
    "F26" = 120
    "G26" = 5532
    "B46" = 39599.03
    "F72" = 50
    "G72" = 2788
    "B85" = 216581.07
    "F168" = 16
    "G168" = 791.6799999999999
    "B180" = 43751.2
    "F182" = 6
    "G182" = 860.64
    "F186" = 61
    "G186" = 1026.02
    "B198" = 63265.91
    "F204" = 20
    "G204" = 14425.2
    "B205" = 39855.5
    "F220" = 70
    "G220" = 5251.4
    "F224" = 86
    "G224" = 6990.08
    "F226" = 76
    "G226" = 3473.2
    "B228" = 46203.28
    "F247" = 57
    "G247" = 1791.51
    "F252" = 82
    "G252" = 4235.3
    "F253" = 78
    "G253" = 3626.22
    "F256" = 12
    "G256" = 557.88
    "F258" = 24
    "G258" = 3042.96
    "F264" = 63
    "G264" = 1513.89
    "F265" = 72
    "G265" = 3767.04
    "B267" = 38218.56
    "F283" = 15
    "G283" = 9404.85
    "F285" = 35
    "G285" = 2750.65
    "B298" = 140124.18
    "F313" = 75
    "G313" = 10165.5
    "F320" = 4
    "G320" = 1183
    "F322" = 42
    "G322" = 1942.5
    "F338" = 123
    "G338" = 4630.95
    "B349" = 162713.38
    "F359" = 31
    "G359" = 4987.9
    "F364" = 59
    "G364" = 2280.35
    "F370" = 38
    "G370" = 2507.24
    "F382" = 87
    "G382" = 9554.34
    "F399" = 384
    "G399" = 22540.8
    "F408" = 7
    "G408" = 3681.65
    "F410" = 60
    "G410" = 9070.200000000001
    "F420" = 124
    "G420" = 11486.12
    "B423" = 312371.82
    "F434" = 6
    "G434" = 1968.72
    "B437" = 52423.8
    "F444" = 64
    "G444" = 7070.08
    "F449" = 43
    "G449" = 2070.45
    "B455" = 26452.19
    "F457" = 29
    "G457" = 616.25
    "B481" = 69252.38
    "F495" = 79
    "G495" = 13113.21
    "F496" = 133
    "G496" = 19985.91
    "B497" = 66561.00999999999
    "F626" = 56
    "G626" = 1200.64
    "B638" = 184395.32
    "F666" = 60
    "G666" = 4156.2
    "B667" = 35824.49
    "F670" = 2
    "G670" = 84.18000000000001
    "F678" = 53
    "G678" = 2488.88
    "F684" = 52
    "G684" = 9012.639999999999
    "B688" = 113293.86
    "F744" = 0
    "G744" = 0
    "F764" = 1
    "G764" = 3215.12
    "F769" = 43
    "G769" = 27149.77
    "B773" = 228629.12
    "F780" = 121
    "G780" = 10345.5
    "B785" = 24195.4
    "F798" = 2
    "G798" = 337.32
    "B805" = 11242.13
    "F814" = 120
    "G814" = 16839.6
    "B815" = 49050.94
    "F817" = 96
    "G817" = 7829.76
    "F819" = 5
    "G819" = 772.6
    "F829" = 132
    "G829" = 7092.36
    "F834" = 4
    "G834" = 576
    "B837" = 259999.14
    "F842" = 80
    "G842" = 13111.2
    "F856" = 554
    "G856" = 16686.48
    "F862" = 52
    "G862" = 2454.92
    "B867" = 286667.08
    "F907" = 375
    "G907" = 61166.25
    "B911" = 62131.58
    "B923" = 4097174.16
    "B924" = 4097174.16
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
